$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '59.124.77'
$ws.Range('E2').Value = '  +4.35%  '

# Row 3
$ws.Range('D3').Value = '3.304.17'
$ws.Range('E3').Value = '  +1.50%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '407.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.15%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.47%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.582'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.17%  '

# Row 8
$ws.Range('E8').Value = '  -0.01%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.627'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.26%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.73'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.52%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0979'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.56%  '

# Row 12
$ws.Range('E12').Value = '  +1.13%  '

# Row 13
$ws.Range('D13').Value = '3.831.61'
$ws.Range('E13').Value = '  +1.69%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.72%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.11%  '

# Row 16
$ws.Range('D16').Value = '3.306.03'
$ws.Range('E16').Value = '  +1.64%  '

# Row 17
$ws.Range('E17').Value = '  -1.29%  '

# Row 18
$ws.Range('D18').Value = '58.879.86'
$ws.Range('E18').Value = '  +4.20%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.62%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.35%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000110'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.91%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.14%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '305.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.65%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.13%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.67%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.68%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.77%  '

# Row 28
$ws.Range('E28').Value = '  +5.43%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.72%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.49'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.17%  '

# Row 31
$ws.Range('E31').Value = '  +3.01%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.41'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.38%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '39.73'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.82%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0511'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.74%  '

# Row 36
$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.11'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.41%  '

# Row 37
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.74'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.04%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.88%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.14%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.71%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '139.06'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.46%  '

# Row 42
$ws.Range('E42').Value = '  +1.94%  '

# Row 43
$ws.Range('E43').Value = '  -1.99%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.04%  '

# Row 45
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.39%  '

# Row 46
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.281'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.27%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +8.52%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.24'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.39%  '

# Row 49
$ws.Range('D49').Value = '2.202.47'
$ws.Range('E49').Value = '  +2.58%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.40'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.02%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.31%  '
